$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1012.6667
$ws.Range("I6").Value = 1012.6667
$ws.Range("K6").Value = 3038.0001
$ws.Range("M6").Value = -2926.0001
$ws.Range("H9").Value = 137.53334
$ws.Range("I9").Value = 94.09999999999999
$ws.Range("J9").Value = 224.4
$ws.Range("K9").Value = 94.09999999999999
$ws.Range("L9").Value = 224.4
$ws.Range("M9").Value = 74.90000000000001
$ws.Range("N9").Value = -562.4
$ws.Range("H12").Value = 108.4
$ws.Range("I12").Value = 85
$ws.Range("J12").Value = 202
$ws.Range("K12").Value = 85
$ws.Range("L12").Value = 202
$ws.Range("M12").Value = 85
$ws.Range("N12").Value = -542
$ws.Range("H48").Value = 1374.8334
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 1449.8
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 4349.4
$ws.Range("M48").Value = -2708
$ws.Range("N48").Value = -4933.4
$ws.Range("H56").Value = 1374.8334
$ws.Range("I56").Value = 1000
$ws.Range("J56").Value = 1449.8
$ws.Range("K56").Value = 3000
$ws.Range("L56").Value = 4349.4
$ws.Range("M56").Value = -2466
$ws.Range("N56").Value = -5417.4
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H132").Value = 6033.3145
$ws.Range("I132").Value = 3346.963
$ws.Range("J132").Value = 15099.75
$ws.Range("K132").Value = 10040.889
$ws.Range("L132").Value = 45299.25
$ws.Range("M132").Value = -7510.889000000001
$ws.Range("N132").Value = -50359.25
$ws.Range("H137").Value = 3822.5
$ws.Range("I137").Value = 2034.9412
$ws.Range("J137").Value = 7199
$ws.Range("K137").Value = 6104.8236
$ws.Range("L137").Value = 21597
$ws.Range("M137").Value = -3554.8236
$ws.Range("N137").Value = -26697
$ws.Range("H138").Value = 5684339.5
$ws.Range("I138").Value = 1388.5172
$ws.Range("J138").Value = 16671378
$ws.Range("K138").Value = 4165.5516
$ws.Range("L138").Value = 50014134
$ws.Range("M138").Value = 974.4484000000002
$ws.Range("N138").Value = -50024414

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 201
$ws.Range("I4").Value = 201
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 201
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -85
$ws.Range("N4").ClearContents()
$ws.Range("H8").Value = 4000
$ws.Range("J8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("N8").Value = -10288
$ws.Range("H32").Value = 6864.4116
$ws.Range("I32").Value = 7520.48
$ws.Range("J32").Value = 5042
$ws.Range("K32").Value = 7520.48
$ws.Range("L32").Value = 5042
$ws.Range("M32").Value = -7233.48
$ws.Range("N32").Value = -5616
$ws.Range("H74").Value = 2650.1738
$ws.Range("I74").Value = 1834.3334
$ws.Range("J74").Value = 3540.182
$ws.Range("K74").Value = 1834.3334
$ws.Range("L74").Value = 3540.182
$ws.Range("M74").Value = -960.3334
$ws.Range("N74").Value = -5288.182
$ws.Range("H77").Value = 2650.1738
$ws.Range("I77").Value = 1834.3334
$ws.Range("J77").Value = 3540.182
$ws.Range("K77").Value = 9171.666999999999
$ws.Range("L77").Value = 17700.91
$ws.Range("M77").Value = -4803.666999999999
$ws.Range("N77").Value = -26436.91
$ws.Range("H132").Value = 10419370
$ws.Range("I132").Value = 16669131
$ws.Range("J132").Value = 3101.2222
$ws.Range("K132").Value = 50007393
$ws.Range("L132").Value = 9303.6666
$ws.Range("M132").Value = -50004863
$ws.Range("N132").Value = -14363.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3527.2222
$ws.Range("I134").Value = 2412.0881
$ws.Range("J134").Value = 6974
$ws.Range("K134").Value = 7236.2643
$ws.Range("L134").Value = 20922
$ws.Range("M134").Value = -4701.2643
$ws.Range("N134").Value = -25992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 11835.667
$ws.Range("I10").Value = 401.2
$ws.Range("J10").Value = 69008
$ws.Range("K10").Value = 401.2
$ws.Range("L10").Value = 69008
$ws.Range("M10").Value = -262.2
$ws.Range("N10").Value = -69286
$ws.Range("H31").Value = 3743.775
$ws.Range("I31").Value = 5131
$ws.Range("J31").Value = 862.61536
$ws.Range("K31").Value = 5131
$ws.Range("L31").Value = 862.61536
$ws.Range("M31").Value = -4836
$ws.Range("N31").Value = -1452.61536
$ws.Range("H34").Value = 3743.775
$ws.Range("I34").Value = 5131
$ws.Range("J34").Value = 862.61536
$ws.Range("K34").Value = 5131
$ws.Range("L34").Value = 862.61536
$ws.Range("M34").Value = -4929
$ws.Range("N34").Value = -1266.61536
$ws.Range("H111").Value = 43600.8
$ws.Range("J111").Value = 43600.8
$ws.Range("L111").Value = 43600.8
$ws.Range("N111").Value = -51780.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 731173.4
$ws.Range("I2").Value = 297.75
$ws.Range("J2").Value = 1262719.2
$ws.Range("K2").Value = 1786.5
$ws.Range("L2").Value = 7576315.199999999
$ws.Range("M2").Value = -1673.5
$ws.Range("N2").Value = -7576541.199999999
$ws.Range("H4").Value = 6250162.5
$ws.Range("I4").Value = 6250162.5
$ws.Range("K4").Value = 18750487.5
$ws.Range("M4").Value = -18750375.5
$ws.Range("H7").Value = 213.28572
$ws.Range("I7").Value = 178.2
$ws.Range("K7").Value = 534.5999999999999
$ws.Range("M7").Value = -422.5999999999999
$ws.Range("H34").Value = 3263
$ws.Range("J34").Value = 4333.3335
$ws.Range("L34").Value = 13000.0005
$ws.Range("N34").Value = -13168.0005
$ws.Range("H46").Value = 1700.6666
$ws.Range("J46").Value = 2451
$ws.Range("L46").Value = 7353
$ws.Range("N46").Value = -7535
$ws.Range("H68").Value = 872.2273
$ws.Range("I68").Value = 690.7646999999999
$ws.Range("J68").Value = 986.4815
$ws.Range("K68").Value = 2072.2941
$ws.Range("L68").Value = 2959.4445
$ws.Range("M68").Value = -1261.2941
$ws.Range("N68").Value = -4581.4445
$ws.Range("H71").Value = 872.2273
$ws.Range("I71").Value = 690.7646999999999
$ws.Range("J71").Value = 986.4815
$ws.Range("K71").Value = 6216.882299999999
$ws.Range("L71").Value = 8878.333500000001
$ws.Range("M71").Value = -2160.882299999999
$ws.Range("N71").Value = -16990.3335
$ws.Range("H80").Value = 2692.1538
$ws.Range("I80").Value = 1499.8334
$ws.Range("K80").Value = 4499.5002
$ws.Range("M80").Value = -3563.5002
$ws.Range("H83").Value = 2692.1538
$ws.Range("I83").Value = 1499.8334
$ws.Range("K83").Value = 13498.5006
$ws.Range("M83").Value = -8818.500599999999
$ws.Range("H92").Value = 509.75
$ws.Range("I92").Value = 293
$ws.Range("J92").Value = 639.8
$ws.Range("K92").Value = 879
$ws.Range("L92").Value = 1919.4
$ws.Range("M92").Value = 369
$ws.Range("N92").Value = -4415.4
$ws.Range("H100").Value = 4686.4287
$ws.Range("J100").Value = 4686.4287
$ws.Range("L100").Value = 14059.2861
$ws.Range("N100").Value = -15681.2861
$ws.Range("H103").Value = 5964.381
$ws.Range("I103").Value = 3331.6667
$ws.Range("J103").Value = 6403.1665
$ws.Range("K103").Value = 9995.000100000001
$ws.Range("L103").Value = 19209.4995
$ws.Range("M103").Value = -9116.000100000001
$ws.Range("N103").Value = -20967.4995
$ws.Range("H107").Value = 741.726
$ws.Range("I107").Value = 309
$ws.Range("J107").Value = 1026.9318
$ws.Range("K107").Value = 927
$ws.Range("L107").Value = 3080.7954
$ws.Range("M107").Value = 993
$ws.Range("N107").Value = -6920.7954
$ws.Range("H115").Value = 3083
$ws.Range("I115").Value = 749
$ws.Range("J115").Value = 4250
$ws.Range("K115").Value = 2247
$ws.Range("L115").Value = 12750
$ws.Range("M115").Value = -1072
$ws.Range("N115").Value = -15100
$ws.Range("H116").Value = 806.7
$ws.Range("I116").Value = 674.1111
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2022.3333
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = 1419.6667
$ws.Range("N116").Value = -12884
$ws.Range("H117").Value = 459
$ws.Range("I117").Value = 459
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1377
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 2065
$ws.Range("N117").ClearContents()
$ws.Range("H118").Value = 1150.9615
$ws.Range("I118").Value = 2036
$ws.Range("J118").Value = 885.45
$ws.Range("K118").Value = 6108
$ws.Range("L118").Value = 2656.35
$ws.Range("M118").Value = -4865
$ws.Range("N118").Value = -5142.35
$ws.Range("H131").Value = 1383.2267
$ws.Range("I131").Value = 1992.1111
$ws.Range("J131").Value = 1040.7291
$ws.Range("K131").Value = 5976.3333
$ws.Range("L131").Value = 3122.1873
$ws.Range("M131").Value = -936.3333000000002
$ws.Range("N131").Value = -13202.1873

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3644.0588
$ws.Range("I132").Value = 2813.0386
$ws.Range("J132").Value = 4508.32
$ws.Range("K132").Value = 8439.1158
$ws.Range("L132").Value = 13524.96
$ws.Range("M132").Value = -5909.1158
$ws.Range("N132").Value = -18584.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7148.732
$ws.Range("I132").Value = 3356.5715
$ws.Range("J132").Value = 15316.462
$ws.Range("K132").Value = 10069.7145
$ws.Range("L132").Value = 45949.386
$ws.Range("M132").Value = -7539.7145
$ws.Range("N132").Value = -51009.386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3988.6155
$ws.Range("I132").Value = 4142.1714
$ws.Range("J132").Value = 2645
$ws.Range("K132").Value = 12426.5142
$ws.Range("L132").Value = 7935
$ws.Range("M132").Value = -9896.514200000001
$ws.Range("N132").Value = -12995
